$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force every written value to be stored as literal text (matches the original
# inlineStr cell type) by prefixing with an apostrophe, then strip the resulting
# quote-prefix formatting so the cell keeps its original (unstyled) appearance.
$ws.Range('D2').Value = "'" + '26.493.50'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'" + '  +1.58%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'" + '1.678.04'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'" + '  +2.33%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = "'" + '  +0.03%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'" + '218.11'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'" + '  +2.00%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'" + '0.5338'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'" + '  +1.86%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('E7').Value = "'" + '  +0.01%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = "'" + '0.2699'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'" + '  +3.93%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'" + '0.06419'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'" + '  +1.96%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'" + '21.98'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'" + '  +6.39%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'" + '0.07797'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'" + '  +1.78%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('E12').Value = "'" + '  +2.40%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'" + '1.674.46'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'" + '  +2.09%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('E14').Value = "'" + '  +1.34%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'" + '0.0₅8339'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'" + '  +1.07%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'" + '65.81'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'" + '  +1.28%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'" + '26.531.71'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'" + '  +1.78%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('E18').Value = "'" + '  -0.02%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = "'" + '4.792'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'" + '  +2.22%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = "'" + '193.82'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'" + '  +2.84%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'" + '10.31'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'" + '  +1.42%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'" + '6.344'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'" + '  +2.96%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('E23').Value = "'" + '  +0.05%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'" + '142.38'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'" + '  -2.29%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'" + '0.1284'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'" + '  +5.62%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = "'" + '7.428'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'" + '  +0.15%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'" + '16.35'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'" + '  +3.46%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = "'" + '1.447'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'" + '  +3.61%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = "'" + '0.06279'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'" + '  +5.17%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = "'" + '1.276'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'" + '  +1.67%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = "'" + '3.616'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'" + '  +5.13%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = "'" + '3.466'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'" + '  +1.81%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = "'" + '1.702'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'" + '  +3.67%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = "'" + '1.011'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'" + '  +2.72%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = "'" + '0.6118'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'" + '  +7.87%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('E36').Value = "'" + '  +1.15%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('E37').Value = "'" + '  +1.06%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = "'" + '6.178'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'" + '  +8.22%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = "'" + '0.01634'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'" + '  +1.16%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = "'" + '1.095.33'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'" + '  +5.89%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'" + '0.8640'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'" + '  +1.72%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = "'" + '0.9998'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'" + '  -0.11%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = "'" + '100.59'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'" + '  +0.39%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'" + '1.823.53'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'" + '  +2.01%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'" + '57.95'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'" + '  +3.87%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('B46').Value = "'" + 'BabyDogeCoin'
$ws.Range('B46').Style = 'Normal'
$ws.Range('C46').Value = "'" + 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('C46').Style = 'Normal'
$ws.Range('D46').Value = "'" + '0.0₈105'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'" + '  -2.51%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('B47').Value = "'" + 'EnergySwap'
$ws.Range('B47').Style = 'Normal'
$ws.Range('C47').Value = "'" + 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('C47').Style = 'Normal'
$ws.Range('D47').Value = "'" + '8.167'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'" + '  +1.46%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = "'" + '1.004'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'" + '  +0.01%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = "'" + '0.05210'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'" + '  +1.03%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = "'" + '6.057'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'" + '  +2.42%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = "'" + '1.474'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'" + '  +6.58%  '
$ws.Range('E51').Style = 'Normal'
